# Apply the commit:
#   "pibix última versão e dataset termosporconcelho sem açores"
# -> the "geoMap_ginasios" dataset (Power-Query backed table) is refreshed to a
#    newer version that no longer contains the Azores cities "Ponta Delgada"
#    and "Ponta do Sol" (4 attribute rows each: ginasio / ginásio / ginasios /
#    ginásios), shrinking the table from 384 to 376 data rows (385 -> 377
#    including the header row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("geoMap_ginasios")

# Cities that were removed from the refreshed dataset ("sem açores" = "without
# the Azores").
$citiesToRemove = @("Ponta Delgada", "Ponta do Sol")

# Walk the used range bottom-up (so that row indices of not-yet-visited rows
# are unaffected by the deletions) and physically delete every row whose
# "Cidade" column (A) matches one of the removed cities.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 2; $r--) {
    $city = $ws.Cells.Item($r, 1).Text
    if ($citiesToRemove -contains $city) {
        $ws.Rows.Item($r).Delete()
    }
}

# The table (ListObject) backing the range resizes itself automatically when
# rows are deleted, but the workbook-level "DadosExternos_1" defined name
# (used by the query table / Power Query refresh machinery) does not, so fix
# it up by hand: $A$1:$C$385 -> $A$1:$C$377.
$newLastRow = $ws.UsedRange.Rows.Count
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*DadosExternos_1*") {
        $n.RefersTo = "=geoMap_ginasios!`$A`$1:`$C`$$newLastRow"
    }
}

# Clear the lingering cell selection (J7:J8) that was left over from editing,
# resetting it back to the top-left cell.
$ws.Range("A1").Select()
